# Auto-generated script to update cryptos price/volume columns (D, E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:E range to Text format first so numeric-looking strings
# (e.g. "41.49", "1.000") are stored as text, matching the original inlineStr cells,
# instead of being auto-converted into numbers by Excel.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.546.38"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.755.28"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "324.29"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "0.4549"
$ws.Range("E7").Value = "  +1.70%  "

$ws.Range("D8").Value = "0.3546"
$ws.Range("E8").Value = "  -1.89%  "

$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("D10").Value = "41.49"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").Value = "6.010"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").Value = "7.158"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").Value = "1.754.15"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").Value = "93.60"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").Value = "0.06393"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "17.10"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").Value = "5.737"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("D23").Value = "27.593.45"

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").Value = "2.072"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "165.17"
$ws.Range("E26").Value = "  +1.89%  "

$ws.Range("D27").Value = "20.13"
$ws.Range("E27").Value = "  -1.55%  "

$ws.Range("D28").Value = "1.963.27"
$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").Value = "2.133"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "125.59"
$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("D31").Value = "1.088"
$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("D32").Value = "0.09216"
$ws.Range("E32").Value = "  +2.03%  "

$ws.Range("E33").Value = "  +0.66%  "

$ws.Range("D34").Value = "5.545"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").Value = "11.73"
$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").Value = "0.02283"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("D37").Value = "0.2093"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "0.06019"
$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("D39").Value = "0.6287"
$ws.Range("E39").Value = "  -1.19%  "

$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("E41").Value = "  -2.24%  "

$ws.Range("D42").Value = "1.387"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "7.823"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").Value = "13.12"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").Value = "3.715"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "0.5863"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "122.27"
$ws.Range("E47").Value = "  +0.71%  "

$ws.Range("D48").Value = "1.935"
$ws.Range("E48").Value = "  -1.10%  "

$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("E50").Value = "  -2.69%  "

$ws.Range("D51").Value = "72.24"
$ws.Range("E51").Value = "  -0.29%  "

# Restore the default (no explicit) style on the range so no stray
# number-format/style metadata is left behind, matching the original workbook styling.
$priceVolRange.Style = "Normal"